# Fix the mistaken lab-report number in the title: "№1" -> "№4"
# (commit message: "change mistake from 1 report to 4")

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "Отчет по лабораторной работе №1",  # FindText
    $true,                               # MatchCase
    $false,                              # MatchWholeWord
    $false,                              # MatchWildcards
    $false,                              # MatchSoundsLike
    $false,                              # MatchAllWordForms
    $true,                                # Forward
    1,                                    # Wrap (wdFindContinue)
    $false,                               # Format
    "Отчет по лабораторной работе №4",   # ReplaceWith
    2                                     # Replace (wdReplaceAll)
)
